$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26 (pushes existing row 26.. down to 27..,
# old dimension A1:R116 becomes A1:R117).
$ws.Rows(26).Insert()

# Populate the newly-inserted row 26 with the new price record.
$ws.Range("A26").Value = 5
$ws.Range("B26").Value = "Macroferia Regional de Talca"
$ws.Range("C26").Value = "Maule"
$ws.Range("D26").Value = 44883
$ws.Range("E26").Value = 7
$ws.Range("F26").Value = 100112022
$ws.Range("G26").Value = "Arveja Verde"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 500
$ws.Range("K26").Value = 20000
$ws.Range("L26").Value = 20000
$ws.Range("M26").Value = 20000
$ws.Range("N26").Value = "$/saco 25 kilos"
$ws.Range("O26").Value = "Región del Maule"
$ws.Range("P26").Value = 800
$ws.Range("Q26").Value = 25
$ws.Range("R26").Value = "Hortaliza"
